# "cut content to reduce time"
#
# The deck had two slides removed to shorten the presentation:
#   - "Pattern"                (originally slide index 3)
#   - "Development Everywhere" (originally slide index 15)
#
# Deleting slides automatically renumbers the surviving slides / updates
# p:sldIdLst, p14:sldIdLst (sections) and relationship ids, so we simply
# remove the two slides in question (from the back first isn't required,
# but we delete the lower-index slide first and recompute the second
# slide's new index after that removal).

$p = $ppt.ActivePresentation

# Slide 3 = "Pattern" (title-only slide with Picture2.png)
$p.Slides.Item(3).Delete()

# After removing slide 3, the slide that used to be at position 15
# ("Development Everywhere") is now at position 14.
$p.Slides.Item(14).Delete()
